$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-strings sst.xml reorder (Chequia/Dinamarca, Mali/Cuba/Libano, Mozambique/Niger)
# --- plus refreshed case-count data (dashboard pull) + updated "datos actualizados" timestamp.
# Only cells whose rendered content actually changes are touched below.

# Row 1: (timestamp)
$ws.Range("A1").Value = "Datos actualizados a 10 de Julio de 2020 a las 18:39"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 3247798
$ws.Range("C4").Value = 27799
$ws.Range("D4").Value = 1437703
$ws.Range("E4").Value = 1673970
$ws.Range("G4").Value = 303
$ws.Range("H4").Value = 136125

# Row 6: India
$ws.Range("B6").Value = 819986
$ws.Range("C6").Value = 25144
$ws.Range("D6").Value = 514457
$ws.Range("E6").Value = 283394
$ws.Range("G6").Value = 512
$ws.Range("H6").Value = 22135

# Row 9: Chile
$ws.Range("B9").Value = 309274
$ws.Range("C9").Value = 3058
$ws.Range("D9").Value = 278053
$ws.Range("E9").Value = 24440
$ws.Range("G9").Value = 99
$ws.Range("H9").Value = 6781

# Row 10: España
$ws.Range("B10").Value = 300988
$ws.Range("C10").Value = 852
$ws.Range("G10").Value = 2
$ws.Range("H10").Value = 28403

# Row 15: Italia
$ws.Range("B15").Value = 242639
$ws.Range("C15").Value = 276
$ws.Range("D15").Value = 194273
$ws.Range("E15").Value = 13428
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = 34938

# Row 19: Alemania
$ws.Range("B19").Value = 199332
$ws.Range("C19").Value = 134
$ws.Range("E19").Value = 6606

# Row 23: Canada
$ws.Range("B23").Value = 107021
$ws.Range("C23").Value = 216
$ws.Range("D23").Value = 70812
$ws.Range("E23").Value = 27450
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = 8759

# Row 42: Singapur
$ws.Range("D42").Value = 41780
$ws.Range("E42").Value = 3808

# Row 45: Republica Dominicana
$ws.Range("B45").Value = 41915
$ws.Range("C45").Value = 1125
$ws.Range("D45").Value = 20830
$ws.Range("E45").Value = 20221
$ws.Range("G45").Value = 22
$ws.Range("H45").Value = 864

# Row 47: Israel
$ws.Range("B47").Value = 36266
$ws.Range("C47").Value = 1441
$ws.Range("D47").Value = 18613
$ws.Range("E47").Value = 17302
$ws.Range("G47").Value = 3
$ws.Range("H47").Value = 351

# Row 62: Argelia
$ws.Range("B62").Value = 18242
$ws.Range("C62").Value = 434
$ws.Range("D62").Value = 13124
$ws.Range("E62").Value = 4122
$ws.Range("G62").Value = 8
$ws.Range("H62").Value = 996

# Row 68: Chequia
$ws.Range("A68").Value = "Chequia"
$ws.Range("B68").Value = 13001
$ws.Range("C68").Value = 82
$ws.Range("D68").Value = 8208
$ws.Range("E68").Value = 4441
$ws.Range("H68").Value = 352

# Row 69: Dinamarca
$ws.Range("A69").Value = "Dinamarca"
$ws.Range("B69").Value = 12946
$ws.Range("C69").Value = 30
$ws.Range("D69").Value = 12077
$ws.Range("E69").Value = 260
$ws.Range("H69").Value = 609

# Row 77: Noruega
$ws.Range("B77").Value = 8968
$ws.Range("C77").Value = 3
$ws.Range("E77").Value = 578

# Row 88: Tayikistan
$ws.Range("B88").Value = 6457
$ws.Range("C88").Value = 47
$ws.Range("D88").Value = 5115
$ws.Range("E88").Value = 1287
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 55

# Row 96: Luxemburgo
$ws.Range("B96").Value = 4777
$ws.Range("C96").Value = 58
$ws.Range("D96").Value = 4086
$ws.Range("E96").Value = 581

# Row 100: Grecia
$ws.Range("B100").Value = 3732
$ws.Range("C100").Value = 60
$ws.Range("E100").Value = 2165

# Row 111: Mali
$ws.Range("A111").Value = "Mali"
$ws.Range("B111").Value = 2404
$ws.Range("C111").Value = 34
$ws.Range("D111").Value = 1650
$ws.Range("E111").Value = 633
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 121

# Row 112: Cuba
$ws.Range("A112").Value = "Cuba"
$ws.Range("B112").Value = 2403
$ws.Range("D112").Value = 2244
$ws.Range("E112").Value = 73
$ws.Range("H112").Value = 86

# Row 113: Libano
$ws.Range("A113").Value = "Libano"
$ws.Range("B113").Value = 2082
$ws.Range("C113").Value = 71
$ws.Range("D113").Value = 1402
$ws.Range("E113").Value = 644
$ws.Range("H113").Value = 36

# Row 114: Sudan del Sur
$ws.Range("A114").Value = "Sudan del Sur"
$ws.Range("B114").Value = 2021
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 333
$ws.Range("E114").Value = 1650
$ws.Range("H114").Value = 38

# Row 115: Estonia
$ws.Range("A115").Value = "Estonia"
$ws.Range("B115").Value = 2013
$ws.Range("C115").Value = 2
$ws.Range("D115").Value = 1894
$ws.Range("E115").Value = 50
$ws.Range("H115").Value = 69

# Row 134: Jordania
$ws.Range("B134").Value = 1173
$ws.Range("C134").Value = 4
$ws.Range("D134").Value = 986

# Row 136: Mozambique
$ws.Range("A136").Value = "Mozambique"
$ws.Range("B136").Value = 1111
$ws.Range("C136").Value = 19
$ws.Range("D136").Value = 344
$ws.Range("E136").Value = 758
$ws.Range("H136").Value = 9

# Row 137: Niger
$ws.Range("A137").Value = "Niger"
$ws.Range("B137").Value = 1097
$ws.Range("D137").Value = 976
$ws.Range("E137").Value = 53
$ws.Range("H137").Value = 68
